# edit.ps1 - applies the "added auth and client code for HH creation" edit
#
# Summary of changes (see unified diff):
#  1. "Return a list of the user's Facebook Friends..." bullet gets expanded
#     to "Return a list of households of the user's Facebook Friends ...,
#     but that the user is not a part of" and picks up the (relocated)
#     "_GoBack" bookmark in the middle of the word "households".
#  2. The "Return notification preferences..." bullet gains a
#     <w:lastRenderedPageBreak/> marker in front of its text.
#  3. The "Receiving push notifications..." bullet loses that same
#     <w:lastRenderedPageBreak/> marker (it moved to #2) and has its
#     trailing "etc" wrapped in <w:proofErr spellStart/spellEnd> markers,
#     which means the run is now split in three.
#  4. The stray "_GoBack" bookmark that used to sit after
#     "The vote, "For" or "Against"" is removed (it moved to #1).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Return a list of ... Facebook Friends ..." bullet
# ---------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("Return a list of the user")
$r.Expand(4)  # wdParagraph - grab the whole paragraph (incl. its mark)

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Return a list of </w:t></w:r>
<w:r><w:t>house</w:t></w:r>
<w:r><w:t xml:space="preserve">holds of </w:t></w:r>
<w:r><w:t>the user\u2019s Facebook Friends who are also using the app via a Get request</w:t></w:r>
<w:r><w:t>, but that the user is not a part of</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$xml1 = $xml1.Replace("\u2019", [char]0x2019)
$r.InsertXML($xml1)

# Relocate the "_GoBack" bookmark to sit right in the middle of the new
# word "households" (i.e. between "house" and "holds of ").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$full = $d.Content.Text
$idx = $full.IndexOf("Return a list of households of ")
$bmPos = $idx + ("Return a list of house").Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 2) "Return notification preferences ..." bullet - gains a
#    lastRenderedPageBreak right before its text
# ---------------------------------------------------------------------
$r2 = $d.Content
[void]$r2.Find.Execute("Return notification preferences for the user")
$r2.Expand(4)

$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:lastRenderedPageBreak/><w:t>Return notification preferences for the user for a given household via a get request, and allow a user to update these preferences via a post request:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r2.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) "Receiving push notifications ..." bullet - loses its
#    lastRenderedPageBreak and splits "etc" off into its own run
#    wrapped with proofErr spellStart/spellEnd markers.
# ---------------------------------------------------------------------
$r3 = $d.Content
[void]$r3.Find.Execute("Receiving push notifications for updates")
$r3.Expand(4)

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'
$xml3 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r>$rPr<w:t xml:space="preserve">Receiving push notifications for updates on requests/votes (vote increase, request approved, request failed, </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r>$rPr<w:t>etc</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>$rPr<w:t>).</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r3.InsertXML($xml3)

Write-Host "Edit applied."
